$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.710.37"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.025.85"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.95"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.607"
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.27"
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0812"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.52"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.326.91"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.96"
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.17"
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.028.56"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.627.54"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("E19").Value = "  -2.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.80"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "222.98"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.36"
$ws.Range("E24").Value = "  -3.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.24"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.04"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.127"
$ws.Range("E28").Value = "  -3.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.89"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  -4.77%  "
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  -2.68%  "
$ws.Range("E33").Value = "  +4.70%  "
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("E35").Value = "  -2.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.33"
$ws.Range("E36").Value = "  +5.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.25"
$ws.Range("E37").Value = "  -4.13%  "
$ws.Range("E38").Value = "  -2.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.530.34"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.71"
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.70"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0914"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.17"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.215.88"
$ws.Range("E51").Value = "  -1.61%  "
